$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0) Clear the stale row-outline grouping left over from the old template.
$ws.Rows("1:8").Ungroup()

# 1) Add a new "Default" NPC config row (row 7), mirroring the existing
#    AttackNpc3 / Enemy rows' stat layout.
$ws.Range("A7").Value = "Default"
$ws.Range("B7").Value = "60"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "0"
$ws.Range("E7").Value = "60"
$ws.Range("F7").Value = 600
$ws.Range("G7").Value = 600
$ws.Range("H7").Value = 600
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = 55000
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10
$ws.Range("P7").Value = 10
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = 10
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 5
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 5
$ws.Range("W7").Value = "Prefabs/Object/Alch_plate"
$ws.Range("X7").Value = 2
$ws.Range("Y7").Value = 20
$ws.Range("Z7").Value = "DropBag_1"

# 2) Rename the "Player" NPC row (row 5) to "Player_0_0"
$ws.Range("A5").Value = "Player_0_0"

# 3) Widen column W (Prefab) to fit the longer default prefab path.
$ws.Columns("W").ColumnWidth = 62.71

# 4) Move the sheet selection to A6 (matches the saved view state).
$ws.Range("A6").Select()
